# Update ligand/receptor/edge expression + specificity values with recomputed TPM-based figures
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 0.2160193333333333
$ws.Range("H2").Value = 0.648058
$ws.Range("I2").Value = 0.02486881244588016
$ws.Range("J2").Value = 0.02486881244588016
$ws.Range("M2").Value = 3.173991666666667
$ws.Range("N2").Value = 9.521975000000001
$ws.Range("O2").Value = 0.1364420332266311
$ws.Range("P2").Value = 0.1364420332266311
$ws.Range("Q2").Value = 0.685643563838889
$ws.Range("R2").Value = 6.170792074550001
$ws.Range("S2").Value = 0.003393151334047639
$ws.Range("T2").Value = 0.003393151334047639
$ws.Range("G3").Value = 0.2160193333333333
$ws.Range("H3").Value = 0.648058
$ws.Range("I3").Value = 0.02486881244588016
$ws.Range("J3").Value = 0.02486881244588016
$ws.Range("O3").Value = 0.5095288789807429
$ws.Range("P3").Value = 0.5095288789807428
$ws.Range("Q3").Value = 2.560466069007556
$ws.Range("R3").Value = 23.044194621068
$ws.Range("S3").Value = 0.01267137812713167
$ws.Range("T3").Value = 0.01267137812713167
$ws.Range("G4").Value = 0.2160193333333333
$ws.Range("H4").Value = 0.648058
$ws.Range("I4").Value = 0.02486881244588016
$ws.Range("J4").Value = 0.02486881244588016
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.4770453333333333
$ws.Range("N4").Value = 1.431136
$ws.Range("O4").Value = 0.02050699625485553
$ws.Range("P4").Value = 0.02050699625485553
$ws.Range("Q4").Value = 0.1030510148764444
$ws.Range("R4").Value = 0.9274591338880001
$ws.Range("S4").Value = 0.0005099846436903689
$ws.Range("T4").Value = 0.000509984643690369
$ws.Range("G5").Value = 0.2160193333333333
$ws.Range("H5").Value = 0.648058
$ws.Range("I5").Value = 0.02486881244588016
$ws.Range("J5").Value = 0.02486881244588016
$ws.Range("M5").Value = 7.276137666666667
$ws.Range("N5").Value = 21.828413
$ws.Range("O5").Value = 0.3127831202907618
$ws.Range("P5").Value = 0.3127831202907618
$ws.Range("Q5").Value = 1.571786407994889
$ws.Range("R5").Value = 14.146077671954
$ws.Range("S5").Value = 0.007778544754748128
$ws.Range("T5").Value = 0.00777854475474813
$ws.Range("G6").Value = 0.2160193333333333
$ws.Range("H6").Value = 0.648058
$ws.Range("I6").Value = 0.02486881244588016
$ws.Range("J6").Value = 0.02486881244588016
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.4824416666666667
$ws.Range("N6").Value = 1.447325
$ws.Range("O6").Value = 0.02073897124700851
$ws.Range("P6").Value = 0.02073897124700851
$ws.Range("Q6").Value = 0.1042167272055556
$ws.Range("R6").Value = 0.9379505448500001
$ws.Range("S6").Value = 0.000515753586262356
$ws.Range("T6").Value = 0.0005157535862623561
$ws.Range("I7").Value = 0.9551554900377276
$ws.Range("J7").Value = 0.9551554900377278
$ws.Range("M7").Value = 3.173991666666667
$ws.Range("N7").Value = 9.521975000000001
$ws.Range("O7").Value = 0.1364420332266311
$ws.Range("P7").Value = 0.1364420332266311
$ws.Range("Q7").Value = 26.33403648183612
$ws.Range("R7").Value = 237.006328336525
$ws.Range("S7").Value = 0.1303233571083268
$ws.Range("T7").Value = 0.1303233571083268
$ws.Range("I8").Value = 0.9551554900377276
$ws.Range("J8").Value = 0.9551554900377278
$ws.Range("O8").Value = 0.5095288789807429
$ws.Range("P8").Value = 0.5095288789807428
$ws.Range("S8").Value = 0.4866793060912255
$ws.Range("T8").Value = 0.4866793060912256
$ws.Range("I9").Value = 0.9551554900377276
$ws.Range("J9").Value = 0.9551554900377278
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.4770453333333333
$ws.Range("N9").Value = 1.431136
$ws.Range("O9").Value = 0.02050699625485553
$ws.Range("P9").Value = 0.02050699625485553
$ws.Range("Q9").Value = 3.957959103491555
$ws.Range("R9").Value = 35.62163193142399
$ws.Range("S9").Value = 0.01958737005700838
$ws.Range("T9").Value = 0.01958737005700838
$ws.Range("I10").Value = 0.9551554900377276
$ws.Range("J10").Value = 0.9551554900377278
$ws.Range("M10").Value = 7.276137666666667
$ws.Range("N10").Value = 21.828413
$ws.Range("O10").Value = 0.3127831202907618
$ws.Range("P10").Value = 0.3127831202907618
$ws.Range("Q10").Value = 60.36880209017411
$ws.Range("R10").Value = 543.319218811567
$ws.Range("S10").Value = 0.2987565145368521
$ws.Range("T10").Value = 0.2987565145368522
$ws.Range("I11").Value = 0.9551554900377276
$ws.Range("J11").Value = 0.9551554900377278
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 0.6666666666666666
$ws.Range("M11").Value = 0.4824416666666667
$ws.Range("N11").Value = 1.447325
$ws.Range("O11").Value = 0.02073897124700851
$ws.Range("P11").Value = 0.02073897124700851
$ws.Range("Q11").Value = 4.002731508019444
$ws.Range("R11").Value = 36.024583572175
$ws.Range("S11").Value = 0.01980894224431476
$ws.Range("T11").Value = 0.01980894224431476
$ws.Range("G12").Value = 0.173491
$ws.Range("H12").Value = 0.520473
$ws.Range("I12").Value = 0.01997281943922393
$ws.Range("J12").Value = 0.01997281943922393
$ws.Range("M12").Value = 3.173991666666667
$ws.Range("N12").Value = 9.521975000000001
$ws.Range("O12").Value = 0.1364420332266311
$ws.Range("P12").Value = 0.1364420332266311
$ws.Range("Q12").Value = 0.5506589882416667
$ws.Range("R12").Value = 4.955930894175
$ws.Range("S12").Value = 0.002725132093556095
$ws.Range("T12").Value = 0.002725132093556096
$ws.Range("G13").Value = 0.173491
$ws.Range("H13").Value = 0.520473
$ws.Range("I13").Value = 0.01997281943922393
$ws.Range("J13").Value = 0.01997281943922393
$ws.Range("O13").Value = 0.5095288789807429
$ws.Range("P13").Value = 0.5095288789807428
$ws.Range("Q13").Value = 2.056379917128667
$ws.Range("R13").Value = 18.507419254158
$ws.Range("S13").Value = 0.01017672829895256
$ws.Range("T13").Value = 0.01017672829895256
$ws.Range("G14").Value = 0.173491
$ws.Range("H14").Value = 0.520473
$ws.Range("I14").Value = 0.01997281943922393
$ws.Range("J14").Value = 0.01997281943922393
$ws.Range("K14").Value = 2
$ws.Range("L14").Value = 0.6666666666666666
$ws.Range("M14").Value = 0.4770453333333333
$ws.Range("N14").Value = 1.431136
$ws.Range("O14").Value = 0.02050699625485553
$ws.Range("P14").Value = 0.02050699625485553
$ws.Range("Q14").Value = 0.08276307192533332
$ws.Range("R14").Value = 0.7448676473279999
$ws.Range("S14").Value = 0.0004095825334390708
$ws.Range("T14").Value = 0.0004095825334390708
$ws.Range("G15").Value = 0.173491
$ws.Range("H15").Value = 0.520473
$ws.Range("I15").Value = 0.01997281943922393
$ws.Range("J15").Value = 0.01997281943922393
$ws.Range("M15").Value = 7.276137666666667
$ws.Range("N15").Value = 21.828413
$ws.Range("O15").Value = 0.3127831202907618
$ws.Range("P15").Value = 0.3127831202907618
$ws.Range("Q15").Value = 1.262344399927667
$ws.Range("R15").Value = 11.361099599349
$ws.Range("S15").Value = 0.006247160785204444
$ws.Range("T15").Value = 0.006247160785204446
$ws.Range("G16").Value = 0.173491
$ws.Range("H16").Value = 0.520473
$ws.Range("I16").Value = 0.01997281943922393
$ws.Range("J16").Value = 0.01997281943922393
$ws.Range("K16").Value = 2
$ws.Range("L16").Value = 0.6666666666666666
$ws.Range("M16").Value = 0.4824416666666667
$ws.Range("N16").Value = 1.447325
$ws.Range("O16").Value = 0.02073897124700851
$ws.Range("P16").Value = 0.02073897124700851
$ws.Range("Q16").Value = 0.08369928719166665
$ws.Range("R16").Value = 0.753293584725
$ws.Range("S16").Value = 0.0004142157280717577
$ws.Range("T16").Value = 0.0004142157280717578
$ws.Range("G17").Value = 0.000025
$ws.Range("H17").Value = 0.00007499999999999999
$ws.Range("I17").Value = 0.000002878077168156263
$ws.Range("J17").Value = 0.000002878077168156264
$ws.Range("M17").Value = 3.173991666666667
$ws.Range("N17").Value = 9.521975000000001
$ws.Range("O17").Value = 0.1364420332266311
$ws.Range("P17").Value = 0.1364420332266311
$ws.Range("Q17").Value = 0.00007934979166666667
$ws.Range("R17").Value = 0.000714148125
$ws.Range("S17").Value = 0.0000003926907006063853
$ws.Range("T17").Value = 0.0000003926907006063853
$ws.Range("G18").Value = 0.000025
$ws.Range("H18").Value = 0.00007499999999999999
$ws.Range("I18").Value = 0.000002878077168156263
$ws.Range("J18").Value = 0.000002878077168156264
$ws.Range("O18").Value = 0.5095288789807429
$ws.Range("P18").Value = 0.5095288789807428
$ws.Range("Q18").Value = 0.0002963237166666667
$ws.Range("R18").Value = 0.00266691345
$ws.Range("S18").Value = 0.000001466463433110732
$ws.Range("T18").Value = 0.000001466463433110732
$ws.Range("G19").Value = 0.000025
$ws.Range("H19").Value = 0.00007499999999999999
$ws.Range("I19").Value = 0.000002878077168156263
$ws.Range("J19").Value = 0.000002878077168156264
$ws.Range("K19").Value = 2
$ws.Range("L19").Value = 0.6666666666666666
$ws.Range("M19").Value = 0.4770453333333333
$ws.Range("N19").Value = 1.431136
$ws.Range("O19").Value = 0.02050699625485553
$ws.Range("P19").Value = 0.02050699625485553
$ws.Range("Q19").Value = 0.00001192613333333333
$ws.Range("R19").Value = 0.0001073352
$ws.Range("S19").Value = 0.00000005902071770856569
$ws.Range("T19").Value = 0.0000000590207177085657
$ws.Range("G20").Value = 0.000025
$ws.Range("H20").Value = 0.00007499999999999999
$ws.Range("I20").Value = 0.000002878077168156263
$ws.Range("J20").Value = 0.000002878077168156264
$ws.Range("M20").Value = 7.276137666666667
$ws.Range("N20").Value = 21.828413
$ws.Range("O20").Value = 0.3127831202907618
$ws.Range("P20").Value = 0.3127831202907618
$ws.Range("Q20").Value = 0.0001819034416666667
$ws.Range("R20").Value = 0.001637130975
$ws.Range("S20").Value = 0.0000009002139570935156
$ws.Range("T20").Value = 0.0000009002139570935158
$ws.Range("G21").Value = 0.000025
$ws.Range("H21").Value = 0.00007499999999999999
$ws.Range("I21").Value = 0.000002878077168156263
$ws.Range("J21").Value = 0.000002878077168156264
$ws.Range("K21").Value = 2
$ws.Range("L21").Value = 0.6666666666666666
$ws.Range("M21").Value = 0.4824416666666667
$ws.Range("N21").Value = 1.447325
$ws.Range("O21").Value = 0.02073897124700851
$ws.Range("P21").Value = 0.02073897124700851
$ws.Range("Q21").Value = 0.00001206104166666667
$ws.Range("R21").Value = 0.000108549375
$ws.Range("S21").Value = 0.00000005968835963706443
$ws.Range("T21").Value = 0.00000005968835963706444
